$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the existing "Never Use Alone" row (row 10) to clarify it is the English line.
$ws.Cells.Item(10, 1).Value = "Never Use Alone (English)"

# 2. Insert a new row right below it for the Spanish line (shifts rows 11+ down by one).
$ws.Rows.Item(11).Insert()

# 3. Fill in the new row 11 with the Spanish "Never Use Alone" entry.
#    Name
$ws.Cells.Item(11, 1).Value = "Never Use Alone (Spanish)"
#    Phone (new number)
$ws.Cells.Item(11, 3).Value = "tel:18009285330"
#    Website (same as the English row)
$ws.Cells.Item(11, 4).Value = "https://neverusealone.com/"
#    Type
$ws.Cells.Item(11, 7).Value = "AOD"
#    Loc
$ws.Cells.Item(11, 8).Value = "Off"
#    Info / description (same as the English row)
$ws.Cells.Item(11, 9).Value = "National anonymous hotline - Will stay on the phone with you if you use alone and will call for help to your address if you experience an overdose - Overdose prevention"
